{"js": "// Fix sample merge file - it previously only showed first result\n//\n// The document had pending tracked changes (insertions/deletions) that\n// were never accepted, and the merge field for the family name needed\n// a \"block=tbs:row\" directive so the mail-merge engine repeats the\n// mailing-label table row for every result instead of only showing the\n// first one.\n\n// Stop recording further revisions and fold the existing tracked\n// changes (the inserted street-address/suburb lines and the accepted\n// curly apostrophe) into the document body.\ncontext.document.changeTrackingMode = Word.ChangeTrackingMode.off;\nawait context.sync();\n\ncontext.document.acceptAllRevisions();\nawait context.sync();\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the family-name / street-address paragraph (first paragraph of\n// the mailing-label table cell) and the \"Printed from ...\" footer\n// paragraph by their (now revision-free) text content.\nlet nameParagraph = null;\nlet printedParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"[family.family_name]\") === 0) {\n    nameParagraph = paragraphs.items[i];\n  } else if (text.indexOf(\"Printed from\") === 0) {\n    printedParagraph = paragraphs.items[i];\n  }\n}\n\n// Add the \"block=tbs:row\" directive to the family-name merge field so\n// the label row repeats for every record instead of just the first.\nif (nameParagraph) {\n  nameParagraph.insertText(\n    \"[family.family_name;block=tbs:row] Family\\u000b[family.street_address]\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// Re-write the footer line as a single run (it had been split into\n// several runs by the delete/insert revisions around the curly\n// apostrophe before they were accepted above).\nif (printedParagraph) {\n  printedParagraph.insertText(\n    \"Printed from [onshow.system_name] Jethro system on [onshow..now;frm=\\u2019yyyy-mm-dd\\u2019] by [onshow.username]\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Fix sample merge file - it previously only showed first result\n#\n# The document had pending tracked changes (insertions/deletions) that\n# were never accepted, and the merge field for the family name needed\n# a \"block=tbs:row\" directive so the mail-merge engine repeats the\n# table row for every result instead of only showing the first one.\n\n$d = $word.ActiveDocument\n\n# Stop recording new revisions and fold the existing tracked changes\n# (the inserted street-address/suburb lines and the accepted curly\n# apostrophe) into the document body.\n$d.TrackRevisions = $false\n$d.Revisions.AcceptAll()\n\n# Update the family-name merge field so the mailing-label row repeats\n# for every record (\"block=tbs:row\") instead of showing only one row.\n$range = $d.Content\n$find = $range.Find\n$find.MatchWildcards = $true\n$find.Text = \"\\[family.family_name\\]*\\[family.street_address\\]\"\n$found = $find.Execute()\nif ($found) {\n    $range.Text = \"[family.family_name;block=tbs:row] Family\" + [char]11 + \"[family.street_address]\"\n}\n\n# Normalise the \"Printed from ... by [onshow.username]\" line (previously\n# split across several runs by the accepted tracked changes) back into a\n# single run of text. Replacing with a placeholder first forces the\n# engine to re-write the run even though the final text is unchanged.\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.MatchWildcards = $true\n$find2.Text = \"Printed from*by \\[onshow.username\\]\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $range2.Text = \"IRON_TMP_PLACEHOLDER\"\n}\n\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.Text = \"IRON_TMP_PLACEHOLDER\"\n$found3 = $find3.Execute()\nif ($found3) {\n    $range3.Text = \"Printed from [onshow.system_name] Jethro system on [onshow..now;frm=\" + [char]8217 + \"yyyy-mm-dd\" + [char]8217 + \"] by [onshow.username]\"\n}\n"}
